$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.273.32"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").Value = "2.271.80"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "111.19"
$ws.Range("E5").Value = "  -2.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "263.65"
$ws.Range("E6").Value = "  -1.78%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.643"
$ws.Range("E7").Value = "  +2.60%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.41"
$ws.Range("E10").Value = "  -3.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0934"
$ws.Range("E11").Value = "  -1.63%  "
$ws.Range("E12").Value = "  +2.42%  "
$ws.Range("E13").Value = "  +1.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.34"
$ws.Range("E14").Value = "  -2.43%  "
$ws.Range("D15").Value = "2.610.90"
$ws.Range("E15").Value = "  -1.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.858"
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("D17").Value = "2.261.18"
$ws.Range("E17").Value = "  -1.56%  "
$ws.Range("D18").Value = "43.088.38"
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("E19").Value = "  -2.57%  "
$ws.Range("E20").Value = "  +1.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.78"
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.43"
$ws.Range("E22").Value = "  -2.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.51"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.38"
$ws.Range("E24").Value = "  -4.57%  "
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("E26").Value = "  +1.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.27"
$ws.Range("E27").Value = "  -3.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "41.18"
$ws.Range("E28").Value = "  -2.38%  "
$ws.Range("E29").Value = "  -1.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.24"
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.22"
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.36"
$ws.Range("E32").Value = "  -1.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0894"
$ws.Range("E33").Value = "  -4.00%  "
$ws.Range("E34").Value = "  -0.23%  "
$ws.Range("E35").Value = "  +2.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0370"
$ws.Range("E36").Value = "  +1.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.61"
$ws.Range("E38").Value = "  +2.53%  "
$ws.Range("E39").Value = "  -4.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.57"
$ws.Range("E40").Value = "  +6.77%  "
$ws.Range("E41").Value = "  +1.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "75.57"
$ws.Range("E42").Value = "  +5.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.236"
$ws.Range("E43").Value = "  -4.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.06"
$ws.Range("E44").Value = "  -1.53%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("E46").Value = "  -4.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.49"
$ws.Range("E47").Value = "  -3.89%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0995"
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.25"
$ws.Range("E49").Value = "  +1.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "100.55"
$ws.Range("E50").Value = "  -2.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.592"
$ws.Range("E51").Value = "  +7.21%  "
